$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the discipline name from " Poluição Ambiental II" to " Poluição Atmosférica"
$ws.Range("B3").Value = " Poluição Atmosférica"
$ws.Range("C3").Value = " Poluição Atmosférica"

# Update "Semestre ideal" value from EA-6 to EA-7
$ws.Range("B9").Value = "EA-7"
$ws.Range("C9").Value = "EA-7"

# Remove rows 23-26 (Requisitos section), which are no longer present
$ws.Range("A23:C26").EntireRow.Delete()
